# Apply the new abituryent (applicant) rows 19-21 to the "qabul" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style used by plain (un-styled) data cells elsewhere in the sheet (e.g. C3),
# used below to strip the transient "Text" number-format style we apply so
# that new cells end up with no explicit style, matching the rest of the data.
$plainStyle = $ws.Range("C3").Style

function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $plainStyle
}

# Row 19 used to be a trailing blank row with a custom (shorter) height;
# AutoFit drops that explicit height so it behaves like a normal data row.
$ws.Rows.Item(19).AutoFit() | Out-Null

# ---- Row 19 ----
$ws.Range("A19").Value = "Aydinova Narine Sergeevna"
$ws.Range("B19").Value = "Yurisprudensiya"
$ws.Range("C19").Value = "Rus tili"
$ws.Range("D19").Value = "Kunduzgi"
$ws.Range("E19").Value = "AD6055389"
Set-TextCell "F19" "60402085220078"
$ws.Range("G19").Value = "Toshkent shahri"
$ws.Range("H19").Value = "Mirzo Ulugʻbek tumani"
Set-TextCell "I19" "998909340132"
Set-TextCell "J19" "+998935617938"
Set-TextCell "K19" "2025-04-25"

# ---- Row 20 ----
$ws.Range("A20").Value = "Xudoyorov Muslimjon Mominjon ogli"
$ws.Range("B20").Value = "Yurisprudensiya"
$ws.Range("C20").Value = "O'zbek tili"
$ws.Range("D20").Value = "Sirtqi"
$ws.Range("E20").Value = "AD7761080"
Set-TextCell "F20" "31103914340034"
$ws.Range("G20").Value = "Fargona viloyati"
$ws.Range("H20").Value = "Margʻilon tumani"
Set-TextCell "I20" "998901669999"
Set-TextCell "J20" "+998916588000"
Set-TextCell "K20" "2025-04-27"

# ---- Row 21 ----
$ws.Range("A21").Value = "Adizov Ismoiljon"
$ws.Range("B21").Value = "Yurisprudensiya"
$ws.Range("C21").Value = "O'zbek tili"
$ws.Range("D21").Value = "Sirtqi"
$ws.Range("E21").Value = "AD4325461"
Set-TextCell "F21" "32804881070096"
$ws.Range("G21").Value = "Toshkent shahri"
$ws.Range("H21").Value = "Mirzo Ulugʻbek tumani"
Set-TextCell "I21" "998903490733"
Set-TextCell "J21" "+998936578677"
Set-TextCell "K21" "2025-04-27"

# Match the saved selection state: rows 19:21 selected with A19 active.
$ws.Range("A19:XFD21").Select() | Out-Null
